$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 120 ---
$ws.Range("A120").Value = 119
$ws.Range("B120").Value = "romania"
$ws.Range("C120").Value = "liga-1"
$ws.Range("D120").Value = "2023-2024"
$ws.Range("E120").Value = 45236.70833333334
$ws.Range("F120").Value = "UTA Arad"
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = "Din. Bucuresti"
$ws.Range("I120").Value = 1
$ws.Range("J120").Value = 1.7
$ws.Range("K120").Value = "30/10/2023 17:12"
$ws.Range("L120").Value = 2.23
$ws.Range("M120").Value = "06/11/2023 16:56"
$ws.Range("N120").Value = 3.75
$ws.Range("O120").Value = "30/10/2023 17:12"
$ws.Range("P120").Value = 3.02
$ws.Range("Q120").Value = "06/11/2023 16:56"
$ws.Range("R120").Value = 5.14
$ws.Range("S120").Value = "30/10/2023 17:12"
$ws.Range("T120").Value = 3.78
$ws.Range("U120").Value = "06/11/2023 16:56"
$ws.Range("V120").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-uta-arad-din-bucuresti/j9VCbndm/"

# --- Row 121 ---
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "romania"
$ws.Range("C121").Value = "liga-1"
$ws.Range("D121").Value = "2023-2024"
$ws.Range("E121").Value = 45236.83333333334
$ws.Range("F121").Value = "FC Hermannstadt"
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = "CFR Cluj"
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 3.57
$ws.Range("K121").Value = "30/10/2023 20:12"
$ws.Range("L121").Value = 3.48
$ws.Range("M121").Value = "06/11/2023 19:57"
$ws.Range("N121").Value = 3.21
$ws.Range("O121").Value = "30/10/2023 20:12"
$ws.Range("P121").Value = 3.15
$ws.Range("Q121").Value = "06/11/2023 19:57"
$ws.Range("R121").Value = 2.16
$ws.Range("S121").Value = "30/10/2023 20:12"
$ws.Range("T121").Value = 2.28
$ws.Range("U121").Value = "06/11/2023 19:59"
$ws.Range("V121").Value = "https://www.betexplorer.com/football/romania/liga-1/fc-hermannstadt-cfr-cluj/6yft48YP/"

# Match the workbook's existing formatting convention: column A (Indice) uses the
# bordered/bold/centered style, column E (data_partida) uses the datetime number
# format. Copy those formats from the row directly above (row 119) onto the two
# new rows, the same way a user extending this generated sheet would in Excel.
$ws.Range("A119").Copy() | Out-Null
$ws.Range("A120:A121").PasteSpecial(-4122) | Out-Null

$ws.Range("E119").Copy() | Out-Null
$ws.Range("E120:E121").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
